# Apply the "LinuxForHealth" re-branding edit to the
# StructureDefinition-match-source-reference workbook.
#
# Changes:
#   Metadata sheet:
#     B2 (URL)       : http://ibm.com/...        -> http://linuxforhealth.org/...
#     B3 (Version)   : 7.0.0                     -> 8.0.0
#     B8 (Date)      : 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
#     B9 (Publisher) : Alvearie Team              -> LinuxForHealth Team
#   Elements sheet:
#     AI2 (Constraint(s) for the root "Extension" row) is cleared out.
#     Q5  (Fixed Value for Extension.url) shares the same canonical URL string
#         as Metadata!B2, so it moves to linuxforhealth.org too.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-source-reference"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-source-reference"
